$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24; existing rows 24.. shift down to 25..
$ws.Rows(24).Insert()

# Populate the newly inserted row 24 with the new weekly data point.
$ws.Range("A24").Value = 11
$ws.Range("B24").Value = "Vega Monumental Concepción"
$ws.Range("C24").Value = "Bíobío"
$ws.Range("D24").Value = 44607
$ws.Range("E24").Value = 8
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100102
$ws.Range("H24").Value = "Cítricos"
$ws.Range("I24").Value = 100102004
$ws.Range("J24").Value = "Mandarina"
$ws.Range("K24").Value = "Murcott"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 220
$ws.Range("N24").Value = 9000
$ws.Range("O24").Value = 9500
$ws.Range("P24").Value = 9227
$ws.Range("Q24").Value = "$/caja 15 kilos granel"
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("S24").Value = 615
$ws.Range("T24").Value = 15
